$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.4946301877498627
$ws.Cells.Item(2, 3).Value = 2.815833806991577
$ws.Cells.Item(2, 4).Value = 0.4332236349582672
$ws.Cells.Item(2, 5).Value = 3.122532367706299
$ws.Cells.Item(2, 6).Value = 0.4092326462268829
$ws.Cells.Item(2, 7).Value = 3.225309371948242
$ws.Cells.Item(2, 8).Value = 0.4020873010158539
$ws.Cells.Item(2, 9).Value = 3.124607563018799
$ws.Cells.Item(2, 10).Value = 0.5248571038246155
$ws.Cells.Item(2, 11).Value = 2.787523746490479
$ws.Cells.Item(2, 12).Value = 0.5152433514595032
$ws.Cells.Item(2, 13).Value = 2.707319021224976
$ws.Cells.Item(2, 14).Value = 0.4846267104148865
$ws.Cells.Item(2, 15).Value = 2.779758453369141
$ws.Cells.Item(2, 16).Value = 0.5758704543113708
$ws.Cells.Item(2, 17).Value = 2.559982776641846
$ws.Cells.Item(2, 18).Value = 0.4107050001621246
$ws.Cells.Item(2, 19).Value = 3.026839256286621
$ws.Cells.Item(2, 20).Value = 0.4933743178844452
$ws.Cells.Item(2, 21).Value = 2.97681450843811
$ws.Cells.Item(2, 22).Value = 0.5273687839508057
$ws.Cells.Item(2, 23).Value = 2.842788934707642
$ws.Cells.Item(2, 24).Value = 0.4992205202579498
$ws.Cells.Item(2, 25).Value = 2.821062564849854
$ws.Cells.Item(2, 26).Value = 0.4420578479766846
$ws.Cells.Item(2, 27).Value = 2.951903104782104
$ws.Cells.Item(2, 28).Value = 0.4807292520999908
$ws.Cells.Item(2, 29).Value = 2.813172578811646
$ws.Cells.Item(2, 30).Value = 0.4832409620285034
$ws.Cells.Item(2, 31).Value = 2.875189781188965
$ws.Cells.Item(2, 32).Value = 43.43062973022461
$ws.Cells.Item(3, 2).Value = 0.8270396590232849
$ws.Cells.Item(3, 3).Value = 0.6224359273910522
$ws.Cells.Item(3, 4).Value = 0.6391823887825012
$ws.Cells.Item(3, 5).Value = 1.184240102767944
$ws.Cells.Item(3, 6).Value = 0.7614758610725403
$ws.Cells.Item(3, 7).Value = 0.8915688991546631
$ws.Cells.Item(3, 8).Value = 0.6129395365715027
$ws.Cells.Item(3, 9).Value = 1.207689881324768
$ws.Cells.Item(3, 10).Value = 0.8287285566329956
$ws.Cells.Item(3, 11).Value = 0.5969237685203552
$ws.Cells.Item(3, 12).Value = 0.8205872178077698
$ws.Cells.Item(3, 13).Value = 0.5960666537284851
$ws.Cells.Item(3, 14).Value = 0.8364801406860352
$ws.Cells.Item(3, 15).Value = 0.6601536870002747
$ws.Cells.Item(3, 16).Value = 0.8239216804504395
$ws.Cells.Item(3, 17).Value = 0.5772038698196411
$ws.Cells.Item(3, 18).Value = 0.7815260887145996
$ws.Cells.Item(3, 19).Value = 0.7541300654411316
$ws.Cells.Item(3, 20).Value = 0.8252641558647156
$ws.Cells.Item(3, 21).Value = 0.6244245767593384
$ws.Cells.Item(3, 22).Value = 0.8031352758407593
$ws.Cells.Item(3, 23).Value = 0.6788083910942078
$ws.Cells.Item(3, 24).Value = 0.8272562026977539
$ws.Cells.Item(3, 25).Value = 0.6161686182022095
$ws.Cells.Item(3, 26).Value = 0.8220162987709045
$ws.Cells.Item(3, 27).Value = 0.6207436919212341
$ws.Cells.Item(3, 28).Value = 0.8243980407714844
$ws.Cells.Item(3, 29).Value = 0.6831640005111694
$ws.Cells.Item(3, 30).Value = 0.8262168765068054
$ws.Cells.Item(3, 31).Value = 0.6058982610702515
$ws.Cells.Item(3, 32).Value = 10.91961765289307
$ws.Cells.Item(4, 2).Value = 0.8360904455184937
$ws.Cells.Item(4, 3).Value = 0.3620176613330841
$ws.Cells.Item(4, 4).Value = 0.6539061069488525
$ws.Cells.Item(4, 5).Value = 0.9343296885490417
$ws.Cells.Item(4, 6).Value = 0.7711328864097595
$ws.Cells.Item(4, 7).Value = 0.6043238639831543
$ws.Cells.Item(4, 8).Value = 0.6305646896362305
$ws.Cells.Item(4, 9).Value = 0.9459807872772217
$ws.Cells.Item(4, 10).Value = 0.8373462557792664
$ws.Cells.Item(4, 11).Value = 0.3622974157333374
$ws.Cells.Item(4, 12).Value = 0.8316299915313721
$ws.Cells.Item(4, 13).Value = 0.369905024766922
$ws.Cells.Item(4, 14).Value = 0.843192458152771
$ws.Cells.Item(4, 15).Value = 0.430520236492157
$ws.Cells.Item(4, 16).Value = 0.8340983986854553
$ws.Cells.Item(4, 17).Value = 0.3622550368309021
$ws.Cells.Item(4, 18).Value = 0.7938246726989746
$ws.Cells.Item(4, 19).Value = 0.4872302114963531
$ws.Cells.Item(4, 20).Value = 0.8349211812019348
$ws.Cells.Item(4, 21).Value = 0.3657464981079102
$ws.Cells.Item(4, 22).Value = 0.8095011115074158
$ws.Cells.Item(4, 23).Value = 0.4549477100372314
$ws.Cells.Item(4, 24).Value = 0.8360471129417419
$ws.Cells.Item(4, 25).Value = 0.3635031580924988
$ws.Cells.Item(4, 26).Value = 0.8358305692672729
$ws.Cells.Item(4, 27).Value = 0.361288994550705
$ws.Cells.Item(4, 28).Value = 0.8396847248077393
$ws.Cells.Item(4, 29).Value = 0.4339023530483246
$ws.Cells.Item(4, 30).Value = 0.8351377248764038
$ws.Cells.Item(4, 31).Value = 0.3640555441379547
$ws.Cells.Item(4, 32).Value = 7.202302932739258
$ws.Cells.Item(5, 2).Value = 0.8402044177055359
$ws.Cells.Item(5, 3).Value = 0.3396600186824799
$ws.Cells.Item(5, 4).Value = 0.6589727997779846
$ws.Cells.Item(5, 5).Value = 0.9045972228050232
$ws.Cells.Item(5, 6).Value = 0.7752468585968018
$ws.Cells.Item(5, 7).Value = 0.5727072358131409
$ws.Cells.Item(5, 8).Value = 0.6358045935630798
$ws.Cells.Item(5, 9).Value = 0.9157621264457703
$ws.Cells.Item(5, 10).Value = 0.8396414518356323
$ws.Cells.Item(5, 11).Value = 0.3409495651721954
$ws.Cells.Item(5, 12).Value = 0.8386453986167908
$ws.Cells.Item(5, 13).Value = 0.344240128993988
$ws.Cells.Item(5, 14).Value = 0.8457041382789612
$ws.Cells.Item(5, 15).Value = 0.4035923182964325
$ws.Cells.Item(5, 16).Value = 0.8373029828071594
$ws.Cells.Item(5, 17).Value = 0.3401212692260742
$ws.Cells.Item(5, 18).Value = 0.7977221608161926
$ws.Cells.Item(5, 19).Value = 0.46119424700737
$ws.Cells.Item(5, 20).Value = 0.8388619422912598
$ws.Cells.Item(5, 21).Value = 0.3428969979286194
$ws.Cells.Item(5, 22).Value = 0.8134418725967407
$ws.Cells.Item(5, 23).Value = 0.4227980375289917
$ws.Cells.Item(5, 24).Value = 0.8395548462867737
$ws.Cells.Item(5, 25).Value = 0.3400019407272339
$ws.Cells.Item(5, 26).Value = 0.8396847248077393
$ws.Cells.Item(5, 27).Value = 0.3402486741542816
$ws.Cells.Item(5, 28).Value = 0.8447514176368713
$ws.Cells.Item(5, 29).Value = 0.4030531942844391
$ws.Cells.Item(5, 30).Value = 0.8393815755844116
$ws.Cells.Item(5, 31).Value = 0.3416381478309631
$ws.Cells.Item(5, 32).Value = 6.813461303710938
$ws.Cells.Item(6, 2).Value = 0.8425428867340088
$ws.Cells.Item(6, 3).Value = 0.3289871513843536
$ws.Cells.Item(6, 4).Value = 0.6635198593139648
$ws.Cells.Item(6, 5).Value = 0.889176070690155
$ws.Cells.Item(6, 6).Value = 0.7769790291786194
$ws.Cells.Item(6, 7).Value = 0.5562561750411987
$ws.Cells.Item(6, 8).Value = 0.6384462118148804
$ws.Cells.Item(6, 9).Value = 0.9029620289802551
$ws.Cells.Item(6, 10).Value = 0.8419799208641052
$ws.Cells.Item(6, 11).Value = 0.3308498859405518
$ws.Cells.Item(6, 12).Value = 0.8397713303565979
$ws.Cells.Item(6, 13).Value = 0.3358559906482697
$ws.Cells.Item(6, 14).Value = 0.8457907438278198
$ws.Cells.Item(6, 15).Value = 0.3835130333900452
$ws.Cells.Item(6, 16).Value = 0.8420665264129639
$ws.Cells.Item(6, 17).Value = 0.3294890820980072
$ws.Cells.Item(6, 18).Value = 0.8023557662963867
$ws.Cells.Item(6, 19).Value = 0.448225349187851
$ws.Cells.Item(6, 20).Value = 0.8397713303565979
$ws.Cells.Item(6, 21).Value = 0.3319608867168427
$ws.Cells.Item(6, 22).Value = 0.8200675845146179
$ws.Cells.Item(6, 23).Value = 0.3985214531421661
$ws.Cells.Item(6, 24).Value = 0.8406807780265808
$ws.Cells.Item(6, 25).Value = 0.3291098177433014
$ws.Cells.Item(6, 26).Value = 0.8415468335151672
$ws.Cells.Item(6, 27).Value = 0.3297704458236694
$ws.Cells.Item(6, 28).Value = 0.8454009890556335
$ws.Cells.Item(6, 29).Value = 0.3842870593070984
$ws.Cells.Item(6, 30).Value = 0.8428027033805847
$ws.Cells.Item(6, 31).Value = 0.3300549983978271
$ws.Cells.Item(6, 32).Value = 6.609016418457031
